$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title " - Matlab" -> " - " + spell-checked "Matlab" run (wrapped in
#    proofErr spellStart/spellEnd, same visible text, split into two runs).
# ---------------------------------------------------------------------------
$titleRange = $d.Paragraphs.Item(1).Range
$found = $titleRange.Find.Execute("Matlab", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="51"/><w:szCs w:val="51"/></w:rPr><w:t>Matlab</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $titleRange.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from after "...inleveren (deadline)." to
#    after "...is als volgt:" (a few paragraphs later).
# ---------------------------------------------------------------------------
$bms = $d.Bookmarks
if ($bms.Exists("_GoBack")) {
    $old = $bms.Item("_GoBack")
    $old.Delete()
}

$target = $d.Content
$found2 = $target.Find.Execute("berekening van het tentamencijfer is als volgt:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $endPos = $target.End
    # Adding a bookmark as a *collapsed* range sitting exactly on a
    # paragraph-end mark is unreliable, so temporarily insert a one-char
    # placeholder, bookmark that, then delete the placeholder again -- the
    # (now zero-length) bookmark stays behind at the right spot.
    $placeholder = $d.Range($endPos, $endPos)
    $placeholder.InsertAfter("X")
    $bmRange = $d.Range($endPos, $endPos + 1)
    $bms.Add("_GoBack", $bmRange)
    $d.Range($endPos, $endPos + 1).Text = ""
}

# ---------------------------------------------------------------------------
# 3) Fix the exam-grade formula: "(...)/9 . k" -> "(...)/(9k)"  i.e. move the
#    "greater-than-k" run into the fraction's denominator.
# ---------------------------------------------------------------------------
$om = $d.OMaths.Item(1)
$omRange = $om.Range
$mathXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math"><w:body><w:p><m:oMathPara><m:oMathParaPr><m:jc m:val="left"/></m:oMathParaPr><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>min</m:t></m:r><m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr></m:ctrlPr></m:dPr><m:e><m:nary><m:naryPr><m:chr m:val="&#8721;"/><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr></m:ctrlPr></m:naryPr><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>wk=1</m:t></m:r></m:sub><m:sup><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>k=4</m:t></m:r></m:sup><m:e><m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr></m:ctrlPr></m:dPr><m:e><m:f><m:fPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr></m:ctrlPr></m:fPr><m:num><m:d><m:dPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr></m:ctrlPr></m:dPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>cijfe</m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>r</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>wk</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>-1</m:t></m:r></m:e></m:d></m:num><m:den><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>9</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>&#8729;k</m:t></m:r></m:den></m:f></m:e></m:d></m:e></m:nary><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>+</m:t></m:r><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>cijfer</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>TT</m:t></m:r></m:sub></m:sSub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>,10</m:t></m:r></m:e></m:d><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:sz w:val="18"/></w:rPr><m:t>=eindcijfer</m:t></m:r></m:oMath></m:oMathPara></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$omRange.InsertXML($mathXml)

# ---------------------------------------------------------------------------
# 4) Header date text: "woensdag 22 november 2017" -> "donderdag 7 december
#    2017" (the cached result of a TIME field).
# ---------------------------------------------------------------------------
for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)
    $header = $section.Headers.Item(1)
    if ($header.Exists) {
        $header.Range.Find.Execute("woensdag 22 november 2017", $true, $false, $false, $false, $false, $true, 1, $false, "donderdag 7 december 2017", 2) | Out-Null
    }
}
